$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-rank: "Somalia" jumps ahead of "Crucero", "Sri Lanka", "Guatemala" and
# "Consejo Danes para los Refugiados" (rows 102-106). Somalia gets brand-new
# totals; the other four countries simply slide down one row, each also
# picking up its own small daily update.
$ws.Range("A102").Value = "Somalia"
$ws.Range("B102").Value = 722
$ws.Range("C102").Value = 51
$ws.Range("D102").Value = 44
$ws.Range("E102").Value = 646
$ws.Range("F102").Value = 2
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 32

$ws.Range("A103").Value = "Crucero"
$ws.Range("B103").Value = 712
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 645
$ws.Range("E103").Value = 54
$ws.Range("F103").Value = 4
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 13

$ws.Range("A104").Value = "Sri Lanka"
$ws.Range("B104").Value = 707
$ws.Range("C104").Value = 5
$ws.Range("D104").Value = 184
$ws.Range("E104").Value = 516
$ws.Range("F104").Value = 1
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 7

$ws.Range("A105").Value = "Guatemala"
$ws.Range("B105").Value = 688
$ws.Range("C105").Value = 44
$ws.Range("D105").Value = 72
$ws.Range("E105").Value = 599
$ws.Range("F105").Value = 5
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 17

$ws.Range("A106").Value = "Consejo Danes para los Refugiados"
$ws.Range("B106").Value = 674
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 75
$ws.Range("E106").Value = 566
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 33

# --- "San Cristobal y Nieves" and "Burundi" swap ranks (rows 198-199); the
# underlying per-country figures are unchanged, only their order.
$ws.Range("A198").Value = "San Cristobal y Nieves"
$ws.Range("B198").Value = 15
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 8
$ws.Range("E198").Value = 7
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Burundi"
$ws.Range("B199").Value = 15
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 7
$ws.Range("E199").Value = 7
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

# --- Plain daily updates for countries whose rank does not change.
# Estados Unidos (row 4)
$ws.Range("B4").Value = 1165868
$ws.Range("C4").Value = 5094
$ws.Range("E4").Value = 924406
$ws.Range("G4").Value = 108
$ws.Range("H4").Value = 67552

# Moldavia (row 59)
$ws.Range("B59").Value = 4121
$ws.Range("C59").Value = 69
$ws.Range("E59").Value = 2614
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 125

# Islandia (row 77)
$ws.Range("B77").Value = 1799
$ws.Range("C77").Value = 1
$ws.Range("D77").Value = 1717
$ws.Range("E77").Value = 72

# Isla de Man (row 129)
$ws.Range("B129").Value = 321
$ws.Range("C129").Value = 1
$ws.Range("E129").Value = 28

# Sierra Leona (row 137)
$ws.Range("B137").Value = 157
$ws.Range("C137").Value = 2
$ws.Range("D137").Value = 29
$ws.Range("E137").Value = 120
